$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "is_prefered" (D) column marks for existing rows (2-62)
# and rewrite B (id) / C (speaker_variant) per the refreshed export
# (ids now derived directly from the variant text, no levenshtein grouping).
$data = @(
    ,@('#mo', 'Mo')
    ,@('#leridaen', 'Leridaen')
    ,@('#prins', 'Prins')
    ,@('#arcelo', 'Arcelo')
    ,@('#valer', 'Valer')
    ,@('#arestip', 'Arestip')
    ,@('#molimpo', 'Molimpo')
    ,@('#graef', 'graef')
    ,@('#celia', 'Celia')
    ,@('#grf', 'Grf')
    ,@('#mol', 'Mol')
    ,@('#alberto', 'Alberto')
    ,@('#lau', 'Lau')
    ,@('#koning', 'Koning')
    ,@('#le', 'Le')
    ,@('#princes', 'Princes')
    ,@('#konin', 'Konin')
    ,@('#pros', 'Pros')
    ,@('#moli', 'Moli')
    ,@('#pri', 'Pri')
    ,@('#molimp', 'Molimp')
    ,@('#leri', 'Leri')
    ,@('#albert', 'Albert')
    ,@('#celi', 'Celi')
    ,@('#2.-solda', '2. Solda')
    ,@('#ar.-en-ga', 'Ar. en Ga')
    ,@('#ce', 'Ce')
    ,@('#theo', 'Theo')
    ,@('#milli', 'Milli')
    ,@('#graef', 'Graef')
    ,@('#leridano', 'Leridano')
    ,@('#theodor', 'Theodor')
    ,@('#gra', 'Gra')
    ,@('#valerio', 'Valerio')
    ,@('#the', 'The')
    ,@('#ler', 'Ler')
    ,@('#ko', 'Ko')
    ,@('#2.-diena', '2. Diena')
    ,@('#rufi', 'Rufi')
    ,@('#arce', 'Arce')
    ,@('#molim', 'Molim')
    ,@('#aristip', 'Aristip')
    ,@('#c', 'C')
    ,@('#pagie', 'Pagie')
    ,@('#gr', 'Gr')
    ,@('#prospero', 'Prospero')
    ,@('#arcel', 'Arcel')
    ,@('#mi', 'Mi')
    ,@('#prospe', 'Prospe')
    ,@('#i', 'I')
    ,@('#galo', 'Galo')
    ,@('#vale', 'Vale')
    ,@('#rufino', 'Rufino')
    ,@('#millido', 'Millido')
    ,@('#mil', 'Mil')
    ,@('#arc.-en-ga', 'Arc. en Ga')
    ,@('#g', 'G')
    ,@('#theodo', 'Theodo')
    ,@('#kon', 'Kon')
    ,@('#laura', 'Laura')
    ,@('#ar.-en', 'Ar. en')
    ,@('#prin', 'Prin')
    ,@('#leridan', 'Leridan')
    ,@('#cel', 'Cel')
)

$url = "https://www.dbnl.org/tekst/rode001hert01_01"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 1).Value = $url
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $ws.Cells.Item($row, 4).Value = ""
    # rows 63-65 are brand new; make sure the trailing flag columns exist
    # (blank) too, matching the layout of the pre-existing rows.
    if ($row -gt 62) {
        $ws.Cells.Item($row, 5).Value = ""
        $ws.Cells.Item($row, 6).Value = ""
        $ws.Cells.Item($row, 7).Value = ""
        $ws.Cells.Item($row, 8).Value = ""
    }
}

